$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before current row 13 (Synthesized beam (bmaj)) to hold the
# new "Processed channel range" entry. This shifts rows 13-21 down to 14-22,
# which automatically recreates the blank-row gap before the "Validation
# metrics" section header and adds a new trailing row for "Cleaning".
$ws.Rows.Item(13).Insert()

# Copy formatting for the new row from the row above (row 3), which has the
# same visual pattern (label / blank input / required-file / code-var / BQF).
$ws.Range("A3:E3").Copy($ws.Range("A13:E13"))
$ws.Range("B13").Value = ""

# New row 13 content (order matters for shared-string table ordering):
$ws.Range("A13").Value = "Processed channel range"
$ws.Range("C13").Value = "slurmOutput/<latest_executed>.sh"

# New D6 cell: "askapsoft"
$ws.Range("D7").Copy($ws.Range("D6"))
$ws.Range("D6").Value = "askapsoft"

# Finish row 13 content
$ws.Range("D13").Value = "chan_range"

# Column width adjustments (values chosen so the resulting stored OOXML
# column width is as close as possible to the target 57.33203125 / 46.1640625)
$ws.Columns.Item(3).ColumnWidth = 56.42
$ws.Columns.Item(4).ColumnWidth = 45.25

# Sheet view changes: drop the frozen top-left cell, update selection
$ws.Range("D14").Select()

# Workbook window position
$excel.Windows.Item(1).Left = 940
$excel.Windows.Item(1).Top = 1660
